# Update handback/handoff timestamps to reflect new report generation run.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 01:14:56"
$wsZhCn.Range("H2").Value = "2016-03-24 01:15:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 01:15:00"
$wsDeDe.Range("H2").Value = "2016-03-24 01:15:29"
